$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "588.40") are preserved as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.745.30"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "3.343.87"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "588.40"
$ws.Range("E5").Value = "  +2.57%  "
$ws.Range("D6").Value = "184.03"
$ws.Range("E6").Value = "  +1.54%  "
$ws.Range("E7").Value = "  +3.47%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "3.343.80"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").Value = "6.82"
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "3.925.66"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("D15").Value = "66.675.35"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "26.77"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.349.43"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0000165"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "428.88"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").Value = "5.54"
$ws.Range("E20").Value = "  -2.21%  "
$ws.Range("D21").Value = "13.24"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").Value = "7.43"
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("D23").Value = "72.27"
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "5.69"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").Value = "3.469.18"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "0.205"
$ws.Range("E28").Value = "  +5.64%  "
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").Value = "9.07"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("D33").Value = "22.55"
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "5.27"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").Value = "6.68"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("D38").Value = "160.72"
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("D39").Value = "1.46"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").Value = "1.83"
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "26.73"
$ws.Range("E41").Value = "  -3.59%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.868.96"
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("D43").Value = "4.36"
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("D44").Value = "0.764"
$ws.Range("E44").Value = "  -4.48%  "
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "39.86"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").Value = "6.03"
$ws.Range("E47").Value = "  -2.77%  "
$ws.Range("D48").Value = "2.35"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("D49").Value = "23.43"
$ws.Range("E49").Value = "  -3.59%  "
$ws.Range("D50").Value = "316.16"
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("E51").Value = "  +0.85%  "

# Restore default cell style on column D so no stray number-format
# styling is left behind (values remain text).
$ws.Range("D2:D51").Style = "Normal"
